# Remove the trailing "Ver no Jupiter ..." and "(c) 2020 ..." footer
# paragraphs (plus the blank paragraph that separated them from the
# "Requisitos" section), while keeping the single blank paragraph that
# precedes the final page-break paragraph.

$d = $word.ActiveDocument

# Locate the "LOB1053: Física III (Requisito fraco)" paragraph - the
# anchor right before the block that needs to go away.
$anchor = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "*LOB1053*") {
        $anchor = $candidate
        break
    }
}

if ($anchor -ne $null) {
    # The three paragraphs that must be deleted, in document order:
    #   1) the blank paragraph right after the anchor
    #   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   3) the "(c) 2020 ... Creative Commons Attribution" copyright line
    $pBlank = $anchor.Next()
    $pJupiter = $pBlank.Next()
    $pCopyright = $pJupiter.Next()

    $deleteRange = $d.Range($pBlank.Range.Start, $pCopyright.Range.End)
    $deleteRange.Delete()
}
